# MasterConfig.xlsx update — add SQLite config row and expand the
# "scenario" sheet with dedicated XLSM / XLSX / HTML test scenarios,
# tested and verified using sqlite.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "config" sheet: keep MySQL config but move it down a row, and add a
# new SQLite config entry above it.
# ---------------------------------------------------------------------
$cfg = $wb.Worksheets.Item("config")

$cfg.Range("A2").Value = "SQLite"
$cfg.Range("E2").Value = "local.db"

$cfg.Range("A3").Value = "MySQL"
$cfg.Range("B3").Value = "localhost"
$cfg.Range("C3").Value = "root"
$cfg.Range("D3").Value = "Password"
$cfg.Range("E3").Value = "localdb"

# Normalise formatting for both rows onto the plain "Normal" style (row 2
# used to carry an extra explicit-font variant on A2; that variant is no
# longer used anywhere once the data moves down to row 3), then drop the
# untouched middle cells of the SQLite row (no hostname/username/password
# needed for a file-based database).
$cfg.Range("A2:E3").Style = "Excel Built-in Normal"
$cfg.Range("B2:D2").Clear()

# ---------------------------------------------------------------------
# "scenario" sheet: replace the single generic "Test Scenario" row with
# three dedicated scenarios (xlsm / xlsx / html), each querying sqlite.
# ---------------------------------------------------------------------
$scn = $wb.Worksheets.Item("scenario")

$scn.Range("F2").Value = ".xlsm"
$scn.Range("A2").Value = "Test XLSM Scenario"
$scn.Range("A3").Value = "Test XLSX Scenario"
$scn.Range("A4").Value = "Test HTML Scenario"
$scn.Range("F4").Value = ".html"

$scn.Range("D2").Value = "base_dir/test_folder_xlsm/expected_file.csv"
$scn.Range("C2").Value = "base_dir/test_folder_xlsm/actual_file.csv"
$scn.Range("C3").Value = "base_dir/test_folder_xlsx/actual_file.csv"
$scn.Range("D3").Value = "base_dir/test_folder_xlsx/expected_file.csv"
$scn.Range("C4").Value = "base_dir/test_folder_html/actual_file.csv"
$scn.Range("D4").Value = "base_dir/test_folder_html/expected_file.csv"

$scn.Range("E2").Value = "base_dir/test_folder_xlsm/result_file.xlsm"
$scn.Range("E3").Value = "base_dir/test_folder_xlsx/result_file.xlsx"
$scn.Range("E4").Value = "base_dir/test_folder_html/result_file.html"

$scn.Range("B2").Value = "Select 'id;name' union select id||';'||name from user"
$scn.Range("B3").Value = "Select 'id;name' union select id||';'||name from user"
$scn.Range("B4").Value = "Select 'id;name' union select id||';'||name from user"

$scn.Range("F3").Value = ".xlsx"

# ---------------------------------------------------------------------
# Selection / active-cell bookkeeping to mirror the saved view state.
# ---------------------------------------------------------------------
$cfg.Range("D4").Select()
$scn.Range("B4").Select()
